# 복수 Active system ppt 불러오기
#
# Each content slide (S101 / S102 / S103) carries a folded-corner callout
# shape whose first line of text is the relative path to the library used
# when loading the system ("Lib/test "). Update that path so it resolves
# from the sample folder: "../sample/Lib/test ".

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $shp = $s.Shapes.Item($j)
        if (-not $shp.HasTextFrame) { continue }

        $tr = $shp.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count
        for ($k = 1; $k -le $paraCount; $k++) {
            $para = $tr.Paragraphs($k)
            # Paragraph text includes a trailing carriage return for every
            # paragraph but the last one in the frame; strip it for the
            # comparison so we only touch the visible characters below.
            $visibleText = $para.Text.TrimEnd("`r")

            if ($visibleText -eq "Lib/test ") {
                $target = $tr.Characters($para.Start, $visibleText.Length)
                $target.Text = "../sample/Lib/test "
            }
        }
    }
}
